$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing row 557,
# pushing all subsequent rows (old 557-609) down by one (new 558-610).
$ws.Rows(557).Insert()

$ws.Range("A557").Value = 3
$ws.Range("B557").Value = "Femacal de La Calera"
$ws.Range("C557").Value = "Coquimbo"
$ws.Range("D557").Value = "2023-07-25"
$ws.Range("E557").Value = 5
$ws.Range("F557").Value = 100112040
$ws.Range("G557").Value = "Cilantro"
$ws.Range("H557").Value = "Sin especificar"
$ws.Range("I557").Value = "Primera"
$ws.Range("J557").Value = 115
$ws.Range("K557").Value = 3000
$ws.Range("L557").Value = 3500
$ws.Range("M557").Value = 3283
$ws.Range("N557").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O557").Value = "Provincia de Quillota"
$ws.Range("P557").Value = 1094
$ws.Range("Q557").Value = 3
$ws.Range("R557").Value = "Hortaliza"
